# Update functions and Data Model (#50)
# Adds a new "Authorship Resource" column (K) to Sheet1 with a header and
# one data value, giving it its own (unstyled-but-distinct) cell format,
# resizing the new column, and updating the view (zoom + selection) to
# reflect where the author ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column K: header + value -----------------------------------------
$ws.Range("K1").Value = "Authorship Resource"
$ws.Range("K2").Value = "Noémi Villars-Amberg, Daniela Subotic"

# Give the new header/value cells their own cell format (font, no border/
# fill, general number format) distinct from the plain default style.
$ws.Range("K1:K2").Font.Name = "Aptos"
$ws.Range("K1:K2").Font.Size = 12

# --- Column sizing ----------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 32.67

# --- View: zoom + active selection moved to the new column ------------------
$ws.Application.ActiveWindow.Zoom = 90
$ws.Range("K3").Select()
